{"js": "// Insert a new paragraph \"Relatando lo siguiente: <<RELATO>>\" right after\n// the \"Presenta diagn\u00f3stico: <<DX>>.\" paragraph, using the same paragraph\n// formatting (left/right indent + both-justified) as its neighbours.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that ends the diagnosis sentence \"Presenta\n// diagn\u00f3stico: <<DX>>.\" \u2014 match on the stable \"<<DX>>\" placeholder so the\n// lookup survives minor wording/whitespace differences.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"<<DX>>\") >= 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"Presenta diagn\u00f3stico: <<DX>>.\" paragraph.');\n}\n\n// insertParagraph inherits the paragraph/run formatting (indentation,\n// justification, language) from the anchor paragraph, matching the\n// \"Presenta diagn\u00f3stico\" / \"Qu\u00e9 por las circunstancias\" paragraphs around it.\ntarget.insertParagraph(\"Relatando lo siguiente: <<RELATO>>\", \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"Relatando lo siguiente: <<RELATO>>\" right after\n# the \"Presenta diagn\u00f3stico: <<DX>>.\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the \"<<DX>>\" placeholder (the end of the\n# \"Presenta diagn\u00f3stico: <<DX>>.\" sentence) via Find, same as a recorded\n# macro would do.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\"<<DX>>\")\nif (-not $found) {\n    throw 'Could not find the \"Presenta diagn\u00f3stico: <<DX>>.\" paragraph.'\n}\n\n$dxParagraph = $findRange.Paragraphs(1)\n\n# Insert a new empty paragraph right after it; InsertParagraphAfter on a\n# collapsed end-of-paragraph range carries over the paragraph formatting\n# (indentation/justification/language) from $dxParagraph, matching its\n# neighbours.\n$insertionPoint = $dxParagraph.Range\n$insertionPoint.Collapse(0)  # wdCollapseEnd\n$insertionPoint.InsertParagraphAfter()\n\n$newParagraph = $dxParagraph.Next()\n$newParagraph.Range.Text = \"Relatando lo siguiente: <<RELATO>>\"\n"}
